$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1 changes
$ws1.Range("A3").Value = "RECURIMENT2"
$ws1.Columns.Item(1).ColumnWidth = 20.6
$ws1.Range("A7").Select()

# Sheet2 changes
$ws2.Range("A3").Value = "recuriment1"
$ws2.Range("A2").Value = "recuriment2"
$ws2.Range("A5").Value = "plkjh"
$ws2.Columns.Item(1).ColumnWidth = 17.5
$ws2.Range("C11").Select()

$wb.Save()
